$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from APR-2021 to July-2021
$ws.Name = "July-2021"

# Fill in previously empty cells for rows 29-31
$ws.Range("C29").Value = "Leave"
$ws.Range("D29").Value = "Personal Leave"
$ws.Range("D30").Value = "Weekly graph dynamic code generation going on"
$ws.Range("D31").Value = "Weekly Graph dynamic code generation  completed sent to deployment"

# Update existing task description text
$ws.Range("D10").Value = "Montly Target screen search four options implemented and tested"
$ws.Range("D11").Value = "Montly Target screen search options  sent to mohan san"
$ws.Range("D16").Value = "As per client requirement previous month report generation implementation going on in a single result"

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B32").Select()
